$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGCL")

# Build the combined font (black color) + vertical-center alignment format in an
# unused scratch cell so that only a single new cellXfs/font entry is produced,
# then copy that format onto the target cells (mirrors how Excel itself would
# consolidate the formatting into one style record).
$tmpl = $ws.Range("D1")
$tmpl.Font.Color = 0
$tmpl.VerticalAlignment = -4108

# Row 24: "hydrogen" becomes "hydrogen combustion turbine" and now derives its
# lifetime from the natural gas combined cycle w/ CCS row (B20) instead of B4.
$ws.Range("A24").Value = "hydrogen combustion turbine"
$ws.Range("B24").Formula = "=B20"

# New row 25: "hydrogen combined cycle", deriving its lifetime from the
# natural gas combined cycle row (B4), same as hydrogen originally did.
$ws.Range("A25").Value = "hydrogen combined cycle"
$ws.Range("B25").Formula = "=B4"

# Apply the scratch-built format to both new/changed label cells.
$tmpl.Copy()
$ws.Range("A24:A25").PasteSpecial(-4122)
$tmpl.Clear()

$ws.Range("B25").Select()
